$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '26.187.15'
$r.ClearFormats()
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  -0.65%  '
$r.ClearFormats()
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '1.671.71'
$r.ClearFormats()
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  -1.33%  '
$r.ClearFormats()
$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  -0.58%  '
$r.ClearFormats()
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '210.99'
$r.ClearFormats()
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  -3.24%  '
$r.ClearFormats()
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '0.5248'
$r.ClearFormats()
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  -3.36%  '
$r.ClearFormats()
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  -0.59%  '
$r.ClearFormats()
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '0.2639'
$r.ClearFormats()
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  -3.29%  '
$r.ClearFormats()
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.06311'
$r.ClearFormats()
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  -1.85%  '
$r.ClearFormats()
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '21.31'
$r.ClearFormats()
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  -2.42%  '
$r.ClearFormats()
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '0.07541'
$r.ClearFormats()
$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  -1.64%  '
$r.ClearFormats()
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '1.674.06'
$r.ClearFormats()
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  -1.36%  '
$r.ClearFormats()
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '4.456'
$r.ClearFormats()
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  -2.22%  '
$r.ClearFormats()
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '0.5590'
$r.ClearFormats()
$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  -3.67%  '
$r.ClearFormats()
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '66.91'
$r.ClearFormats()
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  +1.10%  '
$r.ClearFormats()
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '0.000007987'
$r.ClearFormats()
$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  -4.24%  '
$r.ClearFormats()
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '26.231.49'
$r.ClearFormats()
$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  -0.73%  '
$r.ClearFormats()
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '1.003'
$r.ClearFormats()
$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  -0.60%  '
$r.ClearFormats()
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '4.783'
$r.ClearFormats()
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  -3.08%  '
$r.ClearFormats()
$r = $ws.Range('B20')
$r.NumberFormat = "@"
$r.Value = 'Avalanche'
$r.ClearFormats()
$r = $ws.Range('C20')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$r.ClearFormats()
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '10.40'
$r.ClearFormats()
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  -4.79%  '
$r.ClearFormats()
$r = $ws.Range('B21')
$r.NumberFormat = "@"
$r.Value = 'BitcoinCash'
$r.ClearFormats()
$r = $ws.Range('C21')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r.ClearFormats()
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '186.25'
$r.ClearFormats()
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  -2.31%  '
$r.ClearFormats()
$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  -1.40%  '
$r.ClearFormats()
$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  -0.61%  '
$r.ClearFormats()
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '149.83'
$r.ClearFormats()
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  +0.96%  '
$r.ClearFormats()
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '0.1246'
$r.ClearFormats()
$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  -4.15%  '
$r.ClearFormats()
$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '7.569'
$r.ClearFormats()
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  -4.03%  '
$r.ClearFormats()
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '15.92'
$r.ClearFormats()
$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  +0.65%  '
$r.ClearFormats()
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '0.06253'
$r.ClearFormats()
$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  +0.71%  '
$r.ClearFormats()
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '1.362'
$r.ClearFormats()
$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  -1.99%  '
$r.ClearFormats()
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '1.280'
$r.ClearFormats()
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  -3.56%  '
$r.ClearFormats()
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '3.498'
$r.ClearFormats()
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  -2.72%  '
$r.ClearFormats()
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '3.432'
$r.ClearFormats()
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  -4.13%  '
$r.ClearFormats()
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '1.636'
$r.ClearFormats()
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  -3.25%  '
$r.ClearFormats()
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '0.9995'
$r.ClearFormats()
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  -3.74%  '
$r.ClearFormats()
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '0.6060'
$r.ClearFormats()
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  -1.46%  '
$r.ClearFormats()
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '2.414'
$r.ClearFormats()
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  +0.07%  '
$r.ClearFormats()
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '2.738'
$r.ClearFormats()
$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  -0.65%  '
$r.ClearFormats()
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '6.122'
$r.ClearFormats()
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  +0.42%  '
$r.ClearFormats()
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '0.01619'
$r.ClearFormats()
$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  -2.21%  '
$r.ClearFormats()
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '1.104.37'
$r.ClearFormats()
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  -0.54%  '
$r.ClearFormats()
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '0.8756'
$r.ClearFormats()
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  -0.97%  '
$r.ClearFormats()
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '1.005'
$r.ClearFormats()
$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  -0.97%  '
$r.ClearFormats()
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '100.15'
$r.ClearFormats()
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  -1.05%  '
$r.ClearFormats()
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '1.824.76'
$r.ClearFormats()
$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  -1.15%  '
$r.ClearFormats()
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '0.00000000109'
$r.ClearFormats()
$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  -0.41%  '
$r.ClearFormats()
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '55.71'
$r.ClearFormats()
$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  -3.50%  '
$r.ClearFormats()
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  -0.61%  '
$r.ClearFormats()
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '8.051'
$r.ClearFormats()
$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  -2.16%  '
$r.ClearFormats()
$r = $ws.Range('D49')
$r.NumberFormat = "@"
$r.Value = '0.05227'
$r.ClearFormats()
$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  -1.22%  '
$r.ClearFormats()
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  -1.02%  '
$r.ClearFormats()
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '5.995'
$r.ClearFormats()
$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  -1.66%  '
$r.ClearFormats()
